$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 208, shifting existing rows 208-221 down to 210-223
$ws.Rows("208:209").Insert()

# New row 208 data
$ws.Cells.Item(208, 1).Value = 7
$ws.Cells.Item(208, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(208, 3).Value = "Ñuble"
$ws.Cells.Item(208, 4).Value = 44931
$ws.Cells.Item(208, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
$ws.Cells.Item(208, 5).Value = 16
$ws.Cells.Item(208, 6).Value = 100112045
$ws.Cells.Item(208, 7).Value = "Zapallo"
$ws.Cells.Item(208, 8).Value = "Camote"
$ws.Cells.Item(208, 9).Value = "1a nueva(o)"
$ws.Cells.Item(208, 10).Value = 400
$ws.Cells.Item(208, 11).Value = 600
$ws.Cells.Item(208, 12).Value = 650
$ws.Cells.Item(208, 13).Value = 625
$ws.Cells.Item(208, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(208, 15).Value = "Región del Maule"
$ws.Cells.Item(208, 16).Value = 625
$ws.Cells.Item(208, 17).Value = 1
$ws.Cells.Item(208, 18).Value = "Hortaliza"

# New row 209 data
$ws.Cells.Item(209, 1).Value = 7
$ws.Cells.Item(209, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(209, 3).Value = "Ñuble"
$ws.Cells.Item(209, 4).Value = 44931
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
$ws.Cells.Item(209, 5).Value = 16
$ws.Cells.Item(209, 6).Value = 100112045
$ws.Cells.Item(209, 7).Value = "Zapallo"
$ws.Cells.Item(209, 8).Value = "Camote"
$ws.Cells.Item(209, 9).Value = "2a nueva(o)"
$ws.Cells.Item(209, 10).Value = 300
$ws.Cells.Item(209, 11).Value = 500
$ws.Cells.Item(209, 12).Value = 550
$ws.Cells.Item(209, 13).Value = 525
$ws.Cells.Item(209, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(209, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(209, 16).Value = 525
$ws.Cells.Item(209, 17).Value = 1
$ws.Cells.Item(209, 18).Value = "Hortaliza"
